$wb = $excel.ActiveWorkbook

# --- Sheet: Summary (4 cell updates) ---
$ws = $wb.Worksheets.Item("Summary")
$ws.Cells.Item(3, 2).Value = 0.01
$ws.Cells.Item(6, 2).Value = 131579.9112553819
$ws.Cells.Item(8, 2).Value = 24253065.61257719
$ws.Cells.Item(10, 2).Value = 2491228.976683192

# --- Sheet: Costs and Revenues (28 cell updates) ---
$ws = $wb.Worksheets.Item("Costs and Revenues")
$ws.Cells.Item(2, 3).Value = 62730.46591140758
$ws.Cells.Item(2, 4).Value = 75394.77333896644
$ws.Cells.Item(2, 5).Value = 91976.24205358134
$ws.Cells.Item(2, 6).Value = 91976.24205358134
$ws.Cells.Item(2, 7).Value = 91976.24205358134
$ws.Cells.Item(2, 8).Value = 91976.24205358134
$ws.Cells.Item(2, 10).Value = 91976.24205358134
$ws.Cells.Item(2, 15).Value = 91976.24205358134
$ws.Cells.Item(3, 3).Value = 196825.9098199031
$ws.Cells.Item(3, 4).Value = 38236.46568336456
$ws.Cells.Item(3, 5).Value = 52530.53686621619
$ws.Cells.Item(5, 3).Value = 38339.65294307929
$ws.Cells.Item(5, 4).Value = 39312.96135688073
$ws.Cells.Item(6, 2).Value = -54153.64424660708
$ws.Cells.Item(6, 3).Value = -189659.5376291605
$ws.Cells.Item(6, 4).Value = -18712.55198267718
$ws.Cells.Item(6, 5).Value = 16649.76056833225
$ws.Cells.Item(6, 6).Value = 69180.29743454844
$ws.Cells.Item(6, 7).Value = 69180.29743454844
$ws.Cells.Item(6, 8).Value = 69180.29743454844
$ws.Cells.Item(6, 9).Value = 69180.29743454844
$ws.Cells.Item(6, 10).Value = 69180.29743454844
$ws.Cells.Item(6, 11).Value = 69180.29743454844
$ws.Cells.Item(6, 12).Value = 69180.29743454844
$ws.Cells.Item(6, 13).Value = 69180.29743454844
$ws.Cells.Item(6, 14).Value = 69180.29743454844
$ws.Cells.Item(6, 15).Value = 69180.29743454844
$ws.Cells.Item(6, 16).Value = 69180.29743454844

# --- Sheet: Installed Capacities (2 cell updates) ---
$ws = $wb.Worksheets.Item("Installed Capacities")
$ws.Cells.Item(3, 3).Value = 216.1492175724446
$ws.Cells.Item(3, 4).Value = 260.7963925174648

# --- Sheet: Added Capacities (3 cell updates) ---
$ws = $wb.Worksheets.Item("Added Capacities")
$ws.Cells.Item(3, 3).Value = 216.1492175724445
$ws.Cells.Item(3, 4).Value = 44.64717494502023
$ws.Cells.Item(3, 5).Value = 65.38503947111997

# --- Sheet: PV Dispatch (93 cell updates) ---
$ws = $wb.Worksheets.Item("PV Dispatch")
$ws.Cells.Item(5, 7).Value = 0.8689415781806812
$ws.Cells.Item(5, 8).Value = 8.899047937542903
$ws.Cells.Item(5, 9).Value = 33.49987019281074
$ws.Cells.Item(5, 10).Value = 73.75033027111266
$ws.Cells.Item(5, 11).Value = 110.5326272755009
$ws.Cells.Item(5, 12).Value = 137.125498098748
$ws.Cells.Item(5, 13).Value = 152.5785378897186
$ws.Cells.Item(5, 14).Value = 155.0474181487245
$ws.Cells.Item(5, 15).Value = 146.4068803306903
$ws.Cells.Item(5, 16).Value = 124.9548851193548
$ws.Cells.Item(5, 17).Value = 93.8359148507591
$ws.Cells.Item(5, 18).Value = 54.58365141039226
$ws.Cells.Item(5, 19).Value = 19.80100621279229
$ws.Cells.Item(5, 20).Value = 3.803791758485934
$ws.Cells.Item(5, 21).Value = 0.06951532625445447
$ws.Cells.Item(6, 7).Value = 0.4649247321369563
$ws.Cells.Item(6, 8).Value = 4.490194123533237
$ws.Cells.Item(6, 9).Value = 16.00727696173293
$ws.Cells.Item(6, 10).Value = 43.92519146913236
$ws.Cells.Item(6, 11).Value = 75.07514852230842
$ws.Cells.Item(6, 12).Value = 100.9478020370177
$ws.Cells.Item(6, 13).Value = 117.8013235769823
$ws.Cells.Item(6, 14).Value = 120.9191740832867
$ws.Cells.Item(6, 15).Value = 110.6174208080415
$ws.Cells.Item(6, 16).Value = 88.78023240253862
$ws.Cells.Item(6, 17).Value = 59.34723422857008
$ws.Cells.Item(6, 18).Value = 28.86611626373139
$ws.Cells.Item(6, 19).Value = 8.635772985087758
$ws.Cells.Item(6, 20).Value = 1.873972933481854
$ws.Cells.Item(6, 21).Value = 0.03058715343006293
$ws.Cells.Item(7, 7).Value = 0.3897772775896541
$ws.Cells.Item(7, 8).Value = 3.465474340751655
$ws.Cells.Item(7, 9).Value = 11.72166576605978
$ws.Cells.Item(7, 10).Value = 27.55725352558855
$ws.Cells.Item(7, 11).Value = 45.28503279632526
$ws.Cells.Item(7, 12).Value = 57.94925088819277
$ws.Cells.Item(7, 13).Value = 61.09935997707642
$ws.Cells.Item(7, 14).Value = 59.64655376060593
$ws.Cells.Item(7, 15).Value = 55.0932464723995
$ws.Cells.Item(7, 16).Value = 47.14179000957051
$ws.Cells.Item(7, 17).Value = 32.63853185343913
$ws.Cells.Item(7, 18).Value = 17.52580377234936
$ws.Cells.Item(7, 19).Value = 6.792754919448789
$ws.Cells.Item(7, 20).Value = 1.665412004246704
$ws.Cells.Item(7, 21).Value = 0.02126057877761752
$ws.Cells.Item(8, 7).Value = 1.048427708612923
$ws.Cells.Item(8, 8).Value = 10.7372102708321
$ws.Cells.Item(8, 9).Value = 40.41950923629976
$ws.Cells.Item(8, 10).Value = 88.98399123388617
$ws.Cells.Item(8, 11).Value = 133.3639361394712
$ws.Cells.Item(8, 12).Value = 165.449755626934
$ws.Cells.Item(8, 13).Value = 184.0947318899791
$ws.Cells.Item(8, 14).Value = 187.0735771170756
$ws.Cells.Item(8, 15).Value = 176.6482740895558
$ws.Cells.Item(8, 16).Value = 150.7652150331742
$ws.Cells.Item(8, 17).Value = 113.2183977184739
$ws.Cells.Item(8, 18).Value = 65.8582970511566
$ws.Cells.Item(8, 19).Value = 23.89104641001701
$ws.Cells.Item(8, 20).Value = 4.589492294453073
$ws.Cells.Item(8, 21).Value = 0.08387421668903385
$ws.Cells.Item(9, 7).Value = 0.560958278245113
$ws.Cells.Item(9, 8).Value = 5.417676003051488
$ws.Cells.Item(9, 9).Value = 19.31369510624622
$ws.Cells.Item(9, 10).Value = 52.99825557753501
$ws.Cells.Item(9, 11).Value = 90.58246021995758
$ws.Cells.Item(9, 12).Value = 121.799296335633
$ws.Cells.Item(9, 13).Value = 142.1340339220183
$ws.Cells.Item(9, 15).Value = 133.4662443858011
$ws.Cells.Item(9, 16).Value = 107.1184277115602
$ws.Cells.Item(9, 17).Value = 71.60583214932356
$ws.Cells.Item(9, 18).Value = 34.82862011770975
$ws.Cells.Item(9, 19).Value = 10.41955398407041
$ws.Cells.Item(9, 20).Value = 2.261055516259907
$ws.Cells.Item(9, 21).Value = 0.03690514988454693
$ws.Cells.Item(10, 7).Value = 0.4702885766708382
$ws.Cells.Item(10, 8).Value = 4.181292981673455
$ws.Cells.Item(10, 9).Value = 14.14286010570121
$ws.Cells.Item(10, 10).Value = 33.24940237062826
$ws.Cells.Item(10, 11).Value = 54.63898190775737
$ws.Cells.Item(10, 12).Value = 69.91908529886263
$ws.Cells.Item(10, 13).Value = 73.71987206850237
$ws.Cells.Item(10, 14).Value = 71.96697828272931
$ws.Cells.Item(10, 15).Value = 66.47315263707451
$ws.Cells.Item(10, 16).Value = 56.87926567298936
$ws.Cells.Item(10, 17).Value = 39.38025527013719
$ws.Cells.Item(10, 18).Value = 21.14588454739968
$ws.Cells.Item(10, 19).Value = 8.195847286163604
$ws.Cells.Item(10, 20).Value = 2.009414827593581
$ws.Cells.Item(10, 21).Value = 0.02565210418204575
$ws.Cells.Item(17, 13).Value = 230.2496698278014
$ws.Cells.Item(19, 10).Value = 41.58545896024957
$ws.Cells.Item(19, 12).Value = 87.44870719248058
$ws.Cells.Item(19, 14).Value = 90.01003351907953

# --- Sheet: Unmet Demand (93 cell updates) ---
$ws = $wb.Worksheets.Item("Unmet Demand")
$ws.Cells.Item(5, 7).Value = 414.4337959369544
$ws.Cells.Item(5, 8).Value = 330.5757541782243
$ws.Cells.Item(5, 9).Value = 176.9760193775952
$ws.Cells.Item(5, 10).Value = 107.2955742555736
$ws.Cells.Item(5, 11).Value = 109.5572237694796
$ws.Cells.Item(5, 12).Value = 98.64091687123928
$ws.Cells.Item(5, 13).Value = 77.7676953375541
$ws.Cells.Item(5, 14).Value = 74.3656454478664
$ws.Cells.Item(5, 15).Value = 83.69133109099639
$ws.Cells.Item(5, 16).Value = 106.2781106359148
$ws.Cells.Item(5, 17).Value = 128.4697750236904
$ws.Cells.Item(5, 18).Value = 161.0018864037399
$ws.Cells.Item(5, 19).Value = 189.2190633734531
$ws.Cells.Item(5, 20).Value = 219.2920578056454
$ws.Cells.Item(5, 21).Value = 251.276137581582
$ws.Cells.Item(6, 7).Value = 136.8785924310737
$ws.Cells.Item(6, 8).Value = 107.7452501129632
$ws.Cells.Item(6, 9).Value = 83.51510303826707
$ws.Cells.Item(6, 10).Value = 82.91243519753434
$ws.Cells.Item(6, 11).Value = 62.76629045205057
$ws.Cells.Item(6, 12).Value = 37.60657774285653
$ws.Cells.Item(6, 13).Value = 24.33271034503603
$ws.Cells.Item(6, 14).Value = 10.42253800004659
$ws.Cells.Item(6, 15).Value = 31.97882363640291
$ws.Cells.Item(6, 16).Value = 45.19417501179163
$ws.Cells.Item(6, 17).Value = 80.63453985745144
$ws.Cells.Item(6, 18).Value = 116.8133877002326
$ws.Cells.Item(6, 19).Value = 163.0473981187501
$ws.Cells.Item(6, 20).Value = 198.2907557613397
$ws.Cells.Item(6, 21).Value = 225.9107949275447
$ws.Cells.Item(7, 7).Value = 167.6012020808691
$ws.Cells.Item(7, 8).Value = 158.7616981666879
$ws.Cells.Item(7, 9).Value = 143.7288091611985
$ws.Cells.Item(7, 10).Value = 99.43816791380756
$ws.Cells.Item(7, 11).Value = 83.72822537421013
$ws.Cells.Item(7, 12).Value = 76.93542539304551
$ws.Cells.Item(7, 13).Value = 77.82642397052864
$ws.Cells.Item(7, 14).Value = 68.03899070462725
$ws.Cells.Item(7, 15).Value = 83.36329197944329
$ws.Cells.Item(7, 16).Value = 90.5862140395771
$ws.Cells.Item(7, 17).Value = 118.8664697760067
$ws.Cells.Item(7, 18).Value = 159.7675876048201
$ws.Cells.Item(7, 19).Value = 217.2238431175235
$ws.Cells.Item(7, 20).Value = 226.2801774240348
$ws.Cells.Item(7, 21).Value = 286.2977687777133
$ws.Cells.Item(8, 7).Value = 414.2543098065221
$ws.Cells.Item(8, 8).Value = 328.737591844935
$ws.Cells.Item(8, 9).Value = 170.0563803341062
$ws.Cells.Item(8, 10).Value = 92.06191329280011
$ws.Cells.Item(8, 11).Value = 86.72591490550931
$ws.Cells.Item(8, 12).Value = 70.31665934305323
$ws.Cells.Item(8, 13).Value = 46.25150133729363
$ws.Cells.Item(8, 14).Value = 42.33948647951536
$ws.Cells.Item(8, 15).Value = 53.44993733213093
$ws.Cells.Item(8, 16).Value = 80.46778072209531
$ws.Cells.Item(8, 17).Value = 109.0872921559756
$ws.Cells.Item(8, 18).Value = 149.7272407629755
$ws.Cells.Item(8, 19).Value = 185.1290231762283
$ws.Cells.Item(8, 20).Value = 218.5063572696783
$ws.Cells.Item(8, 21).Value = 251.2617786911475
$ws.Cells.Item(9, 7).Value = 136.7825588849655
$ws.Cells.Item(9, 8).Value = 106.817768233445
$ws.Cells.Item(9, 9).Value = 80.20868489375378
$ws.Cells.Item(9, 10).Value = 73.83937108913167
$ws.Cells.Item(9, 11).Value = 47.25897875440141
$ws.Cells.Item(9, 12).Value = 16.75508344424118
$ws.Cells.Item(9, 13).Value = 0
$ws.Cells.Item(9, 15).Value = 9.130000058643361
$ws.Cells.Item(9, 16).Value = 26.85597970277001
$ws.Cells.Item(9, 17).Value = 68.37594193669796
$ws.Cells.Item(9, 18).Value = 110.8508838462542
$ws.Cells.Item(9, 19).Value = 161.2636171197674
$ws.Cells.Item(9, 20).Value = 197.9036731785617
$ws.Cells.Item(9, 21).Value = 225.9044769310903
$ws.Cells.Item(10, 7).Value = 167.5206907817879
$ws.Cells.Item(10, 8).Value = 158.0458795257661
$ws.Cells.Item(10, 9).Value = 141.3076148215571
$ws.Cells.Item(10, 10).Value = 93.74601906876785
$ws.Cells.Item(10, 11).Value = 74.37427626277801
$ws.Cells.Item(10, 12).Value = 64.96559098237566
$ws.Cells.Item(10, 13).Value = 65.20591187910269
$ws.Cells.Item(10, 14).Value = 55.71856618250388
$ws.Cells.Item(10, 15).Value = 71.98338581476828
$ws.Cells.Item(10, 16).Value = 80.84873837615825
$ws.Cells.Item(10, 17).Value = 112.1247463593087
$ws.Cells.Item(10, 18).Value = 156.1475068297698
$ws.Cells.Item(10, 19).Value = 215.8207507508087
$ws.Cells.Item(10, 20).Value = 225.9361746006879
$ws.Cells.Item(10, 21).Value = 286.2933772523089
$ws.Cells.Item(17, 12).Value = 28.8362588753597
$ws.Cells.Item(17, 13).Value = 0.09656339947136416
$ws.Cells.Item(19, 12).Value = 47.4359690887577
$ws.Cells.Item(19, 14).Value = 37.67551094615366

# --- Sheet: Household Surplus (2 cell updates) ---
$ws = $wb.Worksheets.Item("Household Surplus")
$ws.Cells.Item(3, 2).Value = 168381.7769200939
$ws.Cells.Item(4, 2).Value = 202375.4442256468
